$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cck"
$ws.Range("C2").Value = "Cckbr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.857847
$ws.Range("H2").Value = 2.573541
$ws.Range("I2").Value = 0.1498657082557423
$ws.Range("J2").Value = 0.1498657082557423
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1771103333333333
$ws.Range("N2").Value = 0.531331
$ws.Range("O2").Value = 0.6758517613986885
$ws.Range("P2").Value = 0.6758517613986885
$ws.Range("Q2").Value = 0.151933568119
$ws.Range("R2").Value = 1.367402113071
$ws.Range("S2").Value = 0.1012870028979054
$ws.Range("T2").Value = 0.1012870028979054

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cck"
$ws.Range("C3").Value = "Cckbr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.857847
$ws.Range("H3").Value = 2.573541
$ws.Range("I3").Value = 0.1498657082557423
$ws.Range("J3").Value = 0.1498657082557423
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.08494466666666667
$ws.Range("N3").Value = 0.254834
$ws.Range("O3").Value = 0.3241482386013114
$ws.Range("P3").Value = 0.3241482386013114
$ws.Range("Q3").Value = 0.072869527466
$ws.Range("R3").Value = 0.655825747194
$ws.Range("S3").Value = 0.04857870535783688
$ws.Range("T3").Value = 0.04857870535783689

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cck"
$ws.Range("C4").Value = "Cckbr"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.378641333333333
$ws.Range("H4").Value = 7.135924
$ws.Range("I4").Value = 0.4155481899527343
$ws.Range("J4").Value = 0.4155481899527343
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1771103333333333
$ws.Range("N4").Value = 0.531331
$ws.Range("O4").Value = 0.6758517613986885
$ws.Range("P4").Value = 0.6758517613986885
$ws.Range("Q4").Value = 0.4212819594271112
$ws.Range("R4").Value = 3.791537634844
$ws.Range("S4").Value = 0.2808489761255923
$ws.Range("T4").Value = 0.2808489761255923

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cck"
$ws.Range("C5").Value = "Cckbr"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.378641333333333
$ws.Range("H5").Value = 7.135924
$ws.Range("I5").Value = 0.4155481899527343
$ws.Range("J5").Value = 0.4155481899527343
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.08494466666666667
$ws.Range("N5").Value = 0.254834
$ws.Range("O5").Value = 0.3241482386013114
$ws.Range("P5").Value = 0.3241482386013114
$ws.Range("Q5").Value = 0.2020528951795556
$ws.Range("R5").Value = 1.818476056616
$ws.Range("S5").Value = 0.134699213827142
$ws.Range("T5").Value = 0.134699213827142

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Cck"
$ws.Range("C6").Value = "Cckbr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.389512333333334
$ws.Range("H6").Value = 7.168537000000001
$ws.Range("I6").Value = 0.4174473515916376
$ws.Range("J6").Value = 0.4174473515916375
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1771103333333333
$ws.Range("N6").Value = 0.531331
$ws.Range("O6").Value = 0.6758517613986885
$ws.Range("P6").Value = 0.6758517613986885
$ws.Range("Q6").Value = 0.4232073258607779
$ws.Range("R6").Value = 3.808865932747
$ws.Range("S6").Value = 0.2821325278644259
$ws.Range("T6").Value = 0.2821325278644258

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cck"
$ws.Range("C7").Value = "Cckbr"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.389512333333334
$ws.Range("H7").Value = 7.168537000000001
$ws.Range("I7").Value = 0.4174473515916376
$ws.Range("J7").Value = 0.4174473515916375
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.08494466666666667
$ws.Range("N7").Value = 0.254834
$ws.Range("O7").Value = 0.3241482386013114
$ws.Range("P7").Value = 0.3241482386013114
$ws.Range("Q7").Value = 0.2029763286508889
$ws.Range("R7").Value = 1.826786957858
$ws.Range("S7").Value = 0.1353148237272117
$ws.Range("T7").Value = 0.1353148237272117

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Cck"
$ws.Range("C8").Value = "Cckbr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.09810400000000001
$ws.Range("H8").Value = 0.294312
$ws.Range("I8").Value = 0.0171387501998857
$ws.Range("J8").Value = 0.0171387501998857
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1771103333333333
$ws.Range("N8").Value = 0.531331
$ws.Range("O8").Value = 0.6758517613986885
$ws.Range("P8").Value = 0.6758517613986885
$ws.Range("Q8").Value = 0.01737523214133333
$ws.Range("R8").Value = 0.156377089272
$ws.Range("S8").Value = 0.01158325451076487
$ws.Range("T8").Value = 0.01158325451076487

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Cck"
$ws.Range("C9").Value = "Cckbr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.09810400000000001
$ws.Range("H9").Value = 0.294312
$ws.Range("I9").Value = 0.0171387501998857
$ws.Range("J9").Value = 0.0171387501998857
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.08494466666666667
$ws.Range("N9").Value = 0.254834
$ws.Range("O9").Value = 0.3241482386013114
$ws.Range("P9").Value = 0.3241482386013114
$ws.Range("Q9").Value = 0.008333411578666667
$ws.Range("R9").Value = 0.07500070420800001
$ws.Range("S9").Value = 0.005555495689120823
$ws.Range("T9").Value = 0.005555495689120823

Write-Output "done"